$d = $word.ActiveDocument

function Replace-InCell($tableIndex, $row, $col, $old, $new) {
    $cell = $d.Tables.Item($tableIndex).Cell($row, $col)
    $rng = $cell.Range
    # wdFindStop (0) + wdReplaceOne (1): replace only the single match
    # inside this cell's own Range, so the edit cannot leak into a
    # neighboring cell/run.
    $rng.Find.Execute($old, $true, $false, $false, $false, $false,
                       $true, 0, $false, $new, 1) | Out-Null
}

Replace-InCell 2 1 2 "nbcnb" "Product Designer"
Replace-InCell 2 2 2 "mnvnb" "NA"
Replace-InCell 2 3 2 "mnvnbv " "Pune"

Replace-InCell 3 1 2 "9889" "20L"
Replace-InCell 3 2 2 "8989" "12L"
Replace-InCell 3 3 2 "898989" "NA"

Replace-InCell 10 2 2 "Only Technical : Yes" "Only Technical : No"
Replace-InCell 10 9 2 "879" "2"
Replace-InCell 10 10 2 "iy" "5"
Replace-InCell 10 11 2 "8989" "NA"
